# edit/submit buttons in ui for LessonModal
# Adds three new "Use Case / Expected Behavior" blocks describing how the
# LessonModal duration / start-time / end-time controls should behave,
# mirroring the existing Use Case blocks already present in the sheet
# (bold label rows "Use Case:" / "Expected Behavior:" followed by
# normal-weight descriptive rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-LabelRow {
    param($row, $text)
    $cell = $ws.Range("A$row")
    $cell.Value = $text
    $cell.Font.Bold = $true
}

function Set-TextRow {
    param($row, $text)
    $cell = $ws.Range("A$row")
    $cell.Value = $text
    $cell.Font.Bold = $false
}

# Bold "Use Case:" / "Expected Behavior:" labels for the three new blocks
Set-LabelRow 90 "Use Case:"
Set-LabelRow 92 "Expected Behavior:"
Set-LabelRow 95 "Use Case:"
Set-LabelRow 97 "Expected Behavior:"
Set-LabelRow 100 "Use Case:"
Set-LabelRow 102 "Expected Behavior:"

# Descriptive text rows, entered in this specific sequence so that the new
# shared-string table entries land in the same order as the source edit.
Set-TextRow 91  "User Changes duration option:"
Set-TextRow 93  "start time and end time are adjusted accordingly"
Set-TextRow 96  "User Changes Start Time"
Set-TextRow 101 "User changes end time:"
Set-TextRow 103 "start time is adjusted to appropriate value based on duration option"
Set-TextRow 98  "end time is adjusted to appropriate value based on duration option"

# Match the saved selection/viewport state from the diff
$ws.Range("A98").Select()
